# Label BOM items better.
# Fixes mislabeled Package/Description values in the BOM worksheet:
#  - Rows 2-4: resistor package name "R-W4" -> "R-1/4W"
#  - Row 6:    Package/Description were swapped; corrected to
#              Package "C-P5mm" / Description "Ceramic Capacitor THT"
#  - Row 7:    Package/Description were swapped; corrected to
#              Package "E-P2.5mm 6.3x11.5mm" / Description "Electrolytic Capacitor THT"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "R-1/4W"
$ws.Range("C3").Value = "R-1/4W"
$ws.Range("C4").Value = "R-1/4W"

$ws.Range("C6").Value = "C-P5mm"
$ws.Range("E6").Value = "Ceramic Capacitor THT"

$ws.Range("C7").Value = "E-P2.5mm 6.3x11.5mm"
$ws.Range("E7").Value = "Electrolytic Capacitor THT"

# Update the saved selection to match the authored workbook.
$ws.Range("E6").Select()
